# ch-lab-order 2.0.0 - ballot
# Updates the "Metadata" sheet (URL, Version, Experimental, new duplicated
# Contact row, Copyright, Count) and appends a new concept row
# (RequestForHistopathExam) to the "Concepts" sheet.

$wb = $excel.ActiveWorkbook

# Helper: write a literal string into a cell without Excel's automatic
# type-coercion turning look-alike numbers/booleans ("true", "7", "1", ...)
# into real numeric/boolean cells. Going through a temporary "=" formula
# and then collapsing it to a value with Paste Special (values only) keeps
# the cell's existing style/number-format untouched, which matters because
# a plain `.Value = "true"` would store a boolean and a quote-prefixed
# literal would allocate a brand-new (unused) cell style.
function Set-TextSafe($cell, [string]$text) {
    if ($text -eq "") {
        $cell.ClearContents()
        return
    }
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------
# Sheet "Metadata"
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
Set-TextSafe $meta.Cells.Item(2, 2) "http://fhir.ch/ig/ch-lab-order/CodeSystem/ServiceRequest.categories"

# Version
Set-TextSafe $meta.Cells.Item(3, 2) "2.0.0-ballot"

# Experimental
Set-TextSafe $meta.Cells.Item(7, 2) "true"

# A new row (duplicating the existing "Contact" row) is inserted right
# after row 10, pushing every later row down by one, through row 23. Give
# the brand-new last row (23) the same formatting as row 22 first, then
# shift the values down (bottom row first, so nothing is overwritten
# before it is read), and finally fill row 11 with a copy of row 10.
$meta.Range("A22:B22").Copy()
$meta.Range("A23:B23").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 22; $r -ge 11; $r--) {
    $destRow = $r + 1
    Set-TextSafe $meta.Cells.Item($destRow, 1) $meta.Cells.Item($r, 1).Text
    Set-TextSafe $meta.Cells.Item($destRow, 2) $meta.Cells.Item($r, 2).Text
}
Set-TextSafe $meta.Cells.Item(11, 1) $meta.Cells.Item(10, 1).Text
Set-TextSafe $meta.Cells.Item(11, 2) $meta.Cells.Item(10, 2).Text

# Copyright (now on row 15 after the shift above)
Set-TextSafe $meta.Cells.Item(15, 2) "CC0-1.0"

# Count (now on row 23 after the shift above)
Set-TextSafe $meta.Cells.Item(23, 2) "7"

# ---------------------------------------------------------------------
# Sheet "Concepts"
# ---------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

# Match the formatting of the previous row before filling in values.
$concepts.Range("A7:D7").Copy()
$concepts.Range("A8:D8").PasteSpecial(-4122)  # xlPasteFormats

Set-TextSafe $concepts.Cells.Item(8, 1) "1"
$concepts.Cells.Item(8, 2).Value = "RequestForHistopathExam"
$concepts.Cells.Item(8, 3).Value = "Anforderung von histopathologischen Untersuchungen"
$concepts.Cells.Item(8, 4).Value = "Electronic ordering of histopathologic tests and/or panels"
